# "add arima call via threads" — append newly-collected sensor readings
# (date, value) to the tail of the data log on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("04/27/2021 18:41:16", "17.622"),
    @("04/27/2021 18:43:06", "18.723"),
    @("04/27/2021 18:51:30", "14.694"),
    @("04/27/2021 18:54:44", "13.75"),
    @("04/27/2021 18:57:56", "13.594")
)

$startRow = 265
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i

    # Column A: timestamp text.
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]

    # Column B: numeric-looking reading that must stay stored as text
    # (matches the rest of the column). Force text via a leading
    # apostrophe, then strip the resulting style tweak so the cell keeps
    # the workbook's default formatting.
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = "'" + $newRows[$i][1]
    $cell.Style = "Normal"
}
